$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3494.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2235.3333
$ws.Range("I113").Value = 2350
$ws.Range("J113").Value = 2006
$ws.Range("K113").Value = 2350
$ws.Range("L113").Value = 2006
$ws.Range("M113").Value = 904
$ws.Range("N113").Value = -8514

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 5507.3726
$ws.Range("I132").Value = 4428.222
$ws.Range("J132").Value = 13601
$ws.Range("K132").Value = 13284.666
$ws.Range("L132").Value = 40803
$ws.Range("M132").Value = -10754.666
$ws.Range("N132").Value = -45863

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1558.4062
$ws.Range("I137").Value = 1060.5333
$ws.Range("J137").Value = 1997.7059
$ws.Range("K137").Value = 3181.5999
$ws.Range("L137").Value = 5993.1177
$ws.Range("M137").Value = -631.5999000000002
$ws.Range("N137").Value = -11093.1177

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5657.676
$ws.Range("I32").Value = 5657.676
$ws.Range("K32").Value = 5657.676
$ws.Range("M32").Value = -5370.676

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1427.0667
$ws.Range("I45").Value = 1308.3636
$ws.Range("J45").Value = 1753.5
$ws.Range("K45").Value = 1308.3636
$ws.Range("L45").Value = 1753.5
$ws.Range("M45").Value = -931.3635999999999
$ws.Range("N45").Value = -2507.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 38462470
$ws.Range("I61").Value = 47619788
$ws.Range("K61").Value = 47619788
$ws.Range("M61").Value = -47619576

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1384.1333
$ws.Range("I74").Value = 720
$ws.Range("J74").Value = 2143.1428
$ws.Range("K74").Value = 720
$ws.Range("L74").Value = 2143.1428
$ws.Range("M74").Value = 154
$ws.Range("N74").Value = -3891.1428

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1384.1333
$ws.Range("I77").Value = 720
$ws.Range("J77").Value = 2143.1428
$ws.Range("K77").Value = 3600
$ws.Range("L77").Value = 10715.714
$ws.Range("M77").Value = 768
$ws.Range("N77").Value = -19451.714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 15154234
$ws.Range("I102").Value = 23812254
$ws.Range("J102").Value = 2696.75
$ws.Range("K102").Value = 23812254
$ws.Range("L102").Value = 2696.75
$ws.Range("M102").Value = -23810632
$ws.Range("N102").Value = -5940.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 3215.36
$ws.Range("I132").Value = 2813.5789
$ws.Range("K132").Value = 8440.736699999999
$ws.Range("M132").Value = -5910.736699999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 38462470
$ws.Range("I136").Value = 47619788
$ws.Range("K136").Value = 142859364
$ws.Range("M136").Value = -142856814

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2299.7827
$ws.Range("I20").Value = 2038.1111
$ws.Range("J20").Value = 3241.8
$ws.Range("K20").Value = 2038.1111
$ws.Range("L20").Value = 3241.8
$ws.Range("M20").Value = -1791.1111
$ws.Range("N20").Value = -3735.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3517.5854
$ws.Range("I134").Value = 806.3103599999999
$ws.Range("K134").Value = 2418.93108
$ws.Range("M134").Value = 116.0689200000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H119").Value = 15000
$ws.Range("J119").Value = 15000
$ws.Range("L119").Value = 15000
$ws.Range("N119").Value = -24676

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 231.83333
$ws.Range("I11").Value = 272.75
$ws.Range("J11").Value = 150
$ws.Range("K11").Value = 818.25
$ws.Range("L11").Value = 450
$ws.Range("M11").Value = -678.25
$ws.Range("N11").Value = -730

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 462.44446
$ws.Range("I114").Value = 341
$ws.Range("J114").Value = 614.25
$ws.Range("K114").Value = 1023
$ws.Range("L114").Value = 1842.75
$ws.Range("M114").Value = 2231
$ws.Range("N114").Value = -8350.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 30304714
$ws.Range("J131").Value = 2175.8333
$ws.Range("L131").Value = 6527.499899999999
$ws.Range("N131").Value = -16607.4999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 50003176
$ws.Range("I70").Value = 41669770
$ws.Range("J70").Value = 66670000
$ws.Range("K70").Value = 41669770
$ws.Range("L70").Value = 66670000
$ws.Range("M70").Value = -41669500
$ws.Range("N70").Value = -66670540

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 50003176
$ws.Range("I73").Value = 41669770
$ws.Range("J73").Value = 66670000
$ws.Range("K73").Value = 41669770
$ws.Range("L73").Value = 66670000
$ws.Range("M73").Value = -41668834
$ws.Range("N73").Value = -66671872

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 679.8182
$ws.Range("I107").Value = 1100.6666
$ws.Range("J107").Value = 522
$ws.Range("K107").Value = 1100.6666
$ws.Range("L107").Value = 522
$ws.Range("M107").Value = 819.3334
$ws.Range("N107").Value = -4362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 9260948
$ws.Range("I122").Value = 1694.3889
$ws.Range("J122").Value = 27779456
$ws.Range("K122").Value = 5083.1667
$ws.Range("L122").Value = 83338368
$ws.Range("M122").Value = -2633.1667
$ws.Range("N122").Value = -83343268

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2631.484
$ws.Range("I132").Value = 2677.7778
$ws.Range("J132").Value = 2567.3845
$ws.Range("K132").Value = 8033.3334
$ws.Range("L132").Value = 7702.1535
$ws.Range("M132").Value = -5503.3334
$ws.Range("N132").Value = -12762.1535

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 962
$ws.Range("I61").Value = 956.6429000000001
$ws.Range("J61").Value = 999.5
$ws.Range("K61").Value = 956.6429000000001
$ws.Range("L61").Value = 999.5
$ws.Range("M61").Value = -754.6429000000001
$ws.Range("N61").Value = -1403.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1000
$ws.Range("I93").Value = 1000
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1000
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = 248
$ws.Range("N93").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H113").Value = 962
$ws.Range("I113").Value = 956.6429000000001
$ws.Range("J113").Value = 999.5
$ws.Range("K113").Value = 956.6429000000001
$ws.Range("L113").Value = 999.5
$ws.Range("M113").Value = 1213.3571
$ws.Range("N113").Value = -5339.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2679.2693
$ws.Range("I132").Value = 2440.9285
$ws.Range("J132").Value = 2957.3333
$ws.Range("K132").Value = 7322.7855
$ws.Range("L132").Value = 8871.999899999999
$ws.Range("M132").Value = -4792.7855
$ws.Range("N132").Value = -13931.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1840.4117
$ws.Range("I136").Value = 1244.3636
$ws.Range("K136").Value = 3733.0908
$ws.Range("M136").Value = -1183.0908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 244.8
$ws.Range("I100").Value = 180.25
$ws.Range("J100").Value = 503
$ws.Range("K100").Value = 360.5
$ws.Range("L100").Value = 1006
$ws.Range("M100").Value = 180.5
$ws.Range("N100").Value = -2088

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 15626261
$ws.Range("I122").Value = 19232156
$ws.Range("K122").Value = 57696468
$ws.Range("M122").Value = -57694018

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3710.2307
$ws.Range("I132").Value = 3830
$ws.Range("J132").Value = 3518.6
$ws.Range("K132").Value = 11490
$ws.Range("L132").Value = 10555.8
$ws.Range("M132").Value = -8960
$ws.Range("N132").Value = -15615.8
